# updated legacy GSC export data
$wb = $excel.ActiveWorkbook

# ---- Chart sheet: rolling date window advanced by one day -----------------
# Drop the oldest day (2025-10-21) so every remaining date/measurement row
# shifts up by one, then append the four newly-exported days at the bottom.
$chart = $wb.Worksheets.Item("Chart")
$chart.Rows.Item(2).Delete()

# Helper so new date cells stay plain text (matches the rest of the column)
# instead of Excel auto-converting a "yyyy-MM-dd"-shaped literal into a date
# serial number.
function Set-TextDate($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

Set-TextDate $chart "A88" "2026-01-16"
$chart.Range("B88").Value = 1.0
$chart.Range("C88").Value = 0.0
$chart.Range("D88").Value = 0.0

Set-TextDate $chart "A89" "2026-01-17"
$chart.Range("B89").Value = 1.0
$chart.Range("C89").Value = 0.0
$chart.Range("D89").Value = 0.0

Set-TextDate $chart "A90" "2026-01-18"
$chart.Range("B90").Value = 1.0
$chart.Range("C90").Value = 0.0
$chart.Range("D90").Value = 0.0

Set-TextDate $chart "A91" "2026-01-19"
$chart.Range("B91").Value = 1.0
$chart.Range("C91").Value = 0.0
# Today's impressions haven't been reported yet -- leave blank.
$chart.Range("D91").Value = ""

# ---- Table sheet: refreshed "Not Started" video count ----------------------
$table = $wb.Worksheets.Item("Table")
$table.Range("C2").Value = 1.0
